$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:D9").Value = 0
